$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

foreach ($row in 2,3) {
    $ws.Range("D$row").Value = 0.0848
    $ws.Range("E$row").Value = 0.0109

    $ws.Range("I$row").Value = 0
    $ws.Range("J$row").Value = 0
    $ws.Range("K$row").Value = 55.5
    $ws.Range("L$row").Value = 0.157312925170068
    $ws.Range("M$row").Value = 44.07
    $ws.Range("N$row").Value = 0.05758526068208546
    $ws.Range("O$row").Value = 0.794054054054054
    $ws.Range("P$row").Value = 35
    $ws.Range("Q$row").Value = 0.04573369920292696
    $ws.Range("R$row").Value = 0.6306306306306306
    $ws.Range("S$row").Value = 9.07
    $ws.Range("T$row").Value = 0.2058089403222147
    $ws.Range("U$row").Value = 3704.1
    $ws.Range("V$row").Value = 4.840062720501765
    $ws.Range("W$row").Value = 0.0549777117384844
    $ws.Range("X$row").Value = 0.06241529226279503
    $ws.Range("Y$row").Value = -0.007437580524310636
    $ws.Range("Z$row").Value = -0.2437305699481865
    $ws.Range("AA$row").Value = -0
    $ws.Range("AB$row").Value = 0.03846785654426859
    $ws.Range("AC$row").Value = -0.03846785654426859
    $ws.Range("AD$row").Value = 915.9
    $ws.Range("AE$row").Value = 0
    $ws.Range("AF$row").Value = 915.9
    $ws.Range("AG$row").Value = -2788.2
    $ws.Range("AH$row").Value = 0.5447894361170593
    $ws.Range("AI$row").Value = 0.4671290865507217
    $ws.Range("AJ$row").Value = 1.378318255969153
    $ws.Range("AK$row").Value = 1.599288746128255

    $ws.Range("AN$row").ClearContents()
    $ws.Range("AP$row").ClearContents()
}
